$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57: phone number A57 was stored as a text string "09876543"; correct it
# to the numeric value 9876543.
$ws.Range("A57").Value = 9876543

# Append the new payment row 58: 76442781 (Cash) 2025-08-18T18:06:01
# Column A keeps its phone number as text (matches source data format).
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "76442781"
$ws.Range("A58").Style = "Normal"

$ws.Range("B58").NumberFormat = "@"
$ws.Range("B58").Value = ""
$ws.Range("B58").Style = "Normal"

$ws.Range("C58").Value = "Cash"
$ws.Range("D58").Value = "2025-08-18T18:06:01"
$ws.Range("E58").Value = 120

$ws.Range("F58").NumberFormat = "@"
$ws.Range("F58").Value = ""
$ws.Range("F58").Style = "Normal"

$ws.Range("G58").Value = 120
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
